# s_Posterior_PS_googleSearchIntervals.xlsx edit script
# - Adds a "Search volume" breakdown table (rows 9-10) to the Losses sheet,
#   mirroring the Bar/Bab/Bay/Total header + search-volume row already
#   present on the Data sheet.
# - Updates the current selection on each sheet and makes "Data" the
#   active/selected tab again.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("Data")
$losses = $wb.Worksheets.Item("Losses")

# --- Losses sheet: add the new rows 9 and 10 ---------------------------

# Row 9: blank corner cell + bold, centered column headers + bold "Total"
$losses.Range("I9").Value = ""
$losses.Range("I9").HorizontalAlignment = -4108  # xlCenter

$losses.Range("J9").Value = "Bar"
$losses.Range("J9").Font.Bold = $true
$losses.Range("J9").HorizontalAlignment = -4108  # xlCenter

$losses.Range("K9").Value = "Bab"
$losses.Range("K9").Font.Bold = $true
$losses.Range("K9").HorizontalAlignment = -4108  # xlCenter

$losses.Range("L9").Value = "Bay"
$losses.Range("L9").Font.Bold = $true
$losses.Range("L9").HorizontalAlignment = -4108  # xlCenter

$losses.Range("M9").Value = "Total"
$losses.Range("M9").Font.Bold = $true

# Row 10: row label + search-volume split + total
$losses.Range("I10").Value = "Search volume"
$losses.Range("I10").Font.Bold = $true
$losses.Range("I10").HorizontalAlignment = -4108  # xlCenter

$losses.Range("J10").Value = 0.5
$losses.Range("J10").NumberFormat = "0%"
$losses.Range("J10").HorizontalAlignment = -4108  # xlCenter

$losses.Range("K10").Value = 0.4
$losses.Range("K10").NumberFormat = "0%"
$losses.Range("K10").HorizontalAlignment = -4108  # xlCenter

$losses.Range("L10").Value = 0.1
$losses.Range("L10").NumberFormat = "0%"
$losses.Range("L10").HorizontalAlignment = -4108  # xlCenter

$losses.Range("M10").Value = 1
$losses.Range("M10").NumberFormat = "0%"

# --- Selections / active tab -------------------------------------------

# Losses sheet selection moves to C6
$losses.Range("C6").Select() | Out-Null

# Data sheet selection moves to K21, and Data becomes the active tab
$data.Range("K21").Select() | Out-Null
$data.Activate()
